$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header "address" -> "location" (H1)
$ws.Range("H1").Value = "location"

# Update interview dates from August to December (L2, M2)
$ws.Range("L2").Value = "14-12-2024"
$ws.Range("M2").Value = "22-12-2024"

# Move active selection cell (cosmetic, matches author's saved view state)
$ws.Range("N6").Select()
